$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add start/end times
$ws.Range("B8").Value = 0.9375
$ws.Range("C8").Value = 0.041666666666666664
$ws.Range("B8").NumberFormat = "h:mm"
$ws.Range("C8").NumberFormat = "h:mm"

# Remove date values from A9 and A10 but keep cell formatting
$ws.Range("A9").ClearContents()
$ws.Range("A10").ClearContents()

# Add new D column formulas for row 3 and row 5 (Total Horas Diarias)
$ws.Range("D3").Formula = "=C3-B3"
$ws.Range("D5").Formula = "=C5-B5"

# D2 and D6 remain blank but get the style (numFmtId 20, time format) applied
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("D6").NumberFormat = "h:mm"

# Selection change to D3
$ws.Range("D3").Select()
